$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = "Globo"
$ws.Cells.Item(8, 2).Value = "RJ TV 1"
$ws.Cells.Item(8, 3).Value = "Social"
$ws.Cells.Item(8, 4).Value = "2025-04-01T12:36"
$ws.Cells.Item(8, 5).Value = "Positivo"
$ws.Cells.Item(8, 6).Value = "Oportunidades de trabalho. Em Campos, 366 vagas, entre elas para taifeiro e nutricionista offshore, Garçom e auxiliar de serviços gerais.  "
